$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = 0.005550000000000003
    "E2" = -0.09965
    "F2" = 0.0424
    "I2" = -0.0008679749190641864
    "J2" = -0.0007781605929032928
    "K2" = 334.3
    "L2" = 0.1224183389482936
    "M2" = 705.4
    "N2" = 0.04850143359071501
    "O2" = 2.11008076577924
    "P2" = 545.9
    "Q2" = 0.03753463651427746
    "R2" = 1.632964403230631
    "S2" = 159.5
    "T2" = 0.2261128437765806
    "U2" = 12446.1
    "V2" = 0.8557608344391806
    "W2" = 0.0893035463105897
    "X2" = 0.04513851340361954
    "Y2" = 0.04416503290697017
    "Z2" = 0.5598768335455888
    "AA2" = 0.0004269482398297536
    "AB2" = 0.04176150485444023
    "AC2" = -0.04133455661461048
    "AD2" = 7027.599999999999
    "AE2" = 54.7013295449024
    "AF2" = 7082.301329544902
    "AG2" = -5363.798670455097
    "AH2" = 0.3274870709665192
    "AI2" = 0.479710417561205
    "AJ2" = -0.5842853447807208
    "AK2" = -2.314375040295805
    "AN2" = 820.0233372228704
    "AP2" = -625.8808250239319
    "D3" = 0.0573
    "E3" = 0.0537
    "I3" = -0.004696385791520666
    "J3" = -0.004683654625218351
    "K3" = 149
    "L3" = 0.2952248860709332
    "M3" = 199.1
    "N3" = 0.1309954602276465
    "O3" = 1.336241610738255
    "P3" = 90.5
    "Q3" = 0.05954339101256662
    "R3" = 0.6073825503355704
    "S3" = 108.6
    "T3" = 0.5454545454545454
    "U3" = 2678.2
    "V3" = 1.762089611158629
    "W3" = 0.1544681733360979
    "X3" = 0.04051832343546788
    "Y3" = 0.11394984990063
    "Z3" = -0.1823141431184625
    "AA3" = 0.0008538964796595071
    "AB3" = 0.03921641033452544
    "AC3" = -0.03836251385486594
    "AD3" = 211.9
    "AE3" = 54.7013295449024
    "AF3" = 266.6013295449024
    "AG3" = -2411.598670455097
    "AH3" = 0.1492309718083541
    "AI3" = 0.2123465131188199
    "AJ3" = 2.704499569595957
    "AK3" = 1.69508745635059
    "AN3" = 24.72578763127188
    "AP3" = -281.4000782328001
    "D4" = -0.0462
    "E4" = -0.253
    "F4" = 0.0248
    "K4" = 185.3
    "L4" = 0.08323974664210952
    "M4" = 506.3
    "N4" = 0.03887438574938575
    "O4" = 2.732325957906098
    "P4" = 455.4
    "Q4" = 0.03496621621621621
    "R4" = 2.45763626551538
    "S4" = 50.89999999999998
    "T4" = 0.1005332806636381
    "U4" = 9767.9
    "V4" = 0.7499923218673218
    "W4" = 0.02413891928508155
    "X4" = 0.0497587033717712
    "Y4" = -0.02561978408668965
    "Z4" = 0.2911533129299747
    "AB4" = 0.04430659937435503
    "AC4" = -0.04430659937435503
    "AD4" = 6815.7
    "AF4" = 6815.7
    "AG4" = -2952.2
    "AH4" = 0.3435384607630155
    "AI4" = 0.5045601930679142
    "AJ4" = -0.2931154312039556
    "AK4" = -0.7892949763387963
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

